# Generate Report for Handoff
# ---------------------------------------------------------------
# The 8c88d43a-...md item has finished its cycle and is dropped from the
# report; the remaining 79d2a641-...md item flips from "Handed back: in
# sync with en-US" to "Ready for handoff" with refreshed timestamps.
#
# Concretely, for every sheet (Overview, zh-cn, de-de):
#   - row 2 (the 79d2a641 entry) gets its Status / handoff-datetime cells
#     refreshed
#   - row 3 (the 8c88d43a entry) is removed entirely, along with its
#     hyperlinks
#
# NOTE on hyperlinks: this runtime's Hyperlinks.Delete() removes every
# hyperlink on the worksheet (it is not scoped to the calling Range), and
# an individual Hyperlink.Delete() is a no-op. So for each sheet we snap-
# shot the hyperlinks we want to keep (row 2) before wiping, delete row 3,
# then recreate the kept hyperlinks and restore their "HyperLink" look
# (underline + the workbook's custom blue) by hand.

$wb = $excel.ActiveWorkbook

function Restyle-Hyperlink($ws, $addr) {
    $ws.Range($addr).Font.Underline = 2
    $ws.Range($addr).Font.Color = 15570276
}

# ---------------- Overview sheet ----------------
$ws = $wb.Worksheets.Item("Overview")

$keep = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 2) {
        $keep += , @($h.Range.Address(), $h.Address, $h.TextToDisplay)
    }
}

$ws.Hyperlinks.Item(1).Range.Hyperlinks.Delete()

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-38-20 08:38:52"

$ws.Rows(3).Delete()

foreach ($l in $keep) {
    $addr = $l[0].Replace("$", "")
    $ws.Hyperlinks.Add($ws.Range($addr), $l[1], "", "", $l[2])
    Restyle-Hyperlink $ws $addr
}

# ---------------- zh-cn sheet ----------------
$ws = $wb.Worksheets.Item("zh-cn")

$keep = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 2) {
        $keep += , @($h.Range.Address(), $h.Address, $h.TextToDisplay)
    }
}

$ws.Hyperlinks.Item(1).Range.Hyperlinks.Delete()

$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-20 08:38:49"

$ws.Rows(3).Delete()

foreach ($l in $keep) {
    $addr = $l[0].Replace("$", "")
    $ws.Hyperlinks.Add($ws.Range($addr), $l[1], "", "", $l[2])
    Restyle-Hyperlink $ws $addr
}

# ---------------- de-de sheet ----------------
$ws = $wb.Worksheets.Item("de-de")

$keep = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 2) {
        $keep += , @($h.Range.Address(), $h.Address, $h.TextToDisplay)
    }
}

$ws.Hyperlinks.Item(1).Range.Hyperlinks.Delete()

$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-20 08:38:52"

$ws.Rows(3).Delete()

foreach ($l in $keep) {
    $addr = $l[0].Replace("$", "")
    $ws.Hyperlinks.Add($ws.Range($addr), $l[1], "", "", $l[2])
    Restyle-Hyperlink $ws $addr
}
